$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.779.40"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "2.508.16"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "2.507.42"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  +4.33%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.54%  "
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").Value = "2.950.86"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").Value = "69.565.93"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "2.472.45"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "2.590.58"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").Value = "0.0₃0895"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "461.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("E46").Value = "  -6.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("E51").Value = "  -0.34%  "
